$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column G header and widen it
$ws.Range("G1").Value = "Role"
$ws.Columns.Item(7).ColumnWidth = 10

# Update row 2 values
$ws.Range("A2").Value = "fyzidejiwi@mailinator.com"
$ws.Range("B2").Value = "fyzidejiwi@mailinator.com"
$ws.Range("C2").Value = "Reiciendis id dolor"
$ws.Range("D2").Value = "Tenetur blanditiis i"
$ws.Range("E2").Value = "+1 (411) 862-1585"
$ws.Range("F2").Value = "Fuga Et in distinct"
$ws.Range("G2").Value = "admin"
